$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = -0.02843597334849807
$ws.Range("C2").Value = 0.5120978290967556
$ws.Range("D2").Value = 0.6508982077648369
$ws.Range("E2").Value = 0.8067826273320694
$ws.Range("F2").Value = 0.8283753329562445
$ws.Range("G2").Value = 19

$ws.Range("B3").Value = -0.05889918563118277
$ws.Range("C3").Value = 0.5375210360923134
$ws.Range("D3").Value = 0.6208141610631079
$ws.Range("E3").Value = 0.7879176105806418
$ws.Range("F3").Value = 0.8084920911945981
$ws.Range("G3").Value = 18

$ws.Range("B4").Value = 0.07850679787661623
$ws.Range("C4").Value = 0.4200275287208594
$ws.Range("D4").Value = 0.2679724280602394
$ws.Range("E4").Value = 0.5176605336127522
$ws.Range("F4").Value = 0.5274203069366139
$ws.Range("G4").Value = 17

$ws.Range("B5").Value = 0.09045959314081749
$ws.Range("C5").Value = 0.4284335755769144
$ws.Range("D5").Value = 0.3144213805930146
$ws.Range("E5").Value = 0.560732895943349
$ws.Range("F5").Value = 0.571536821889252
$ws.Range("G5").Value = 16

$ws.Range("B6").Value = 0.1328729943200521
$ws.Range("C6").Value = 0.4043956606847117
$ws.Range("D6").Value = 0.2696880151153669
$ws.Range("E6").Value = 0.5193149479028761
$ws.Range("F6").Value = 0.5196490393550559
$ws.Range("G6").Value = 15

$ws.Range("B7").Value = 0.1892581571789815
$ws.Range("C7").Value = 0.3380498083070491
$ws.Range("D7").Value = 0.2385894900054021
$ws.Range("E7").Value = 0.4884562314121933
$ws.Range("F7").Value = 0.4672992583619074
$ws.Range("G7").Value = 14

$ws.Range("B8").Value = 0.2086948502624945
$ws.Range("C8").Value = 0.3022036763864308
$ws.Range("D8").Value = 0.1923865449580635
$ws.Range("E8").Value = 0.4386189062934513
$ws.Range("F8").Value = 0.4015417223668171
$ws.Range("G8").Value = 13

$ws.Range("B9").Value = 0.2394124698171911
$ws.Range("C9").Value = 0.2672986829739259
$ws.Range("D9").Value = 0.1375220876090255
$ws.Range("E9").Value = 0.3708397060847524
$ws.Range("F9").Value = 0.2957955502247975
$ws.Range("G9").Value = 12

$ws.Range("B10").Value = 0.195371401722039
$ws.Range("C10").Value = 0.2718328173942094
$ws.Range("D10").Value = 0.1201690274152207
$ws.Range("E10").Value = 0.3466540457216975
$ws.Range("F10").Value = 0.3003313954364828
$ws.Range("G10").Value = 11

$ws.Range("B11").Value = 0.1436448583919225
$ws.Range("C11").Value = 0.2781300272604909
$ws.Range("D11").Value = 0.2246355840812312
$ws.Range("E11").Value = 0.4739573652568669
$ws.Range("F11").Value = 0.4760972574996225
$ws.Range("G11").Value = 10
